$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new staff record in row 2 (Designation, unique_id, and profile links)
$ws.Range("B2").Value = "Professor & Head"
$ws.Range("J2").Value = "VEC-003-01-173"
$ws.Range("D2").Value = "https://scholar.google.co.in/citations?user=M5UgVCMAAAAJ&hl=en"

# Research Gate profile as an actual hyperlink (also sets E2's text)
$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.researchgate.net/profile/Abdul-Razak-Mohamed-Sikkander")

$ws.Range("F2").Value = "https://orcid.org/0000-0002-8458-7448"
$ws.Range("G2").Value = "https://publons.com/researcher/3826750/abdulrazak-mohamed-sikkander/"
$ws.Range("H2").Value = "https://www.scopus.com/authid/detail.uri?authorId=57195150706"
$ws.Range("I2").Value = "https://www.linkedin.com/in/dr-abdul-razak-mohamed-sikkander-51485246?utm_source=share&utm_campaign=share_via&utm_content=profile&utm_medium=android_app"

$ws.Range("L2").Select() | Out-Null
